$d = $word.ActiveDocument

# 1. "The above object has been converted to a JSON string which can be used to store/transmit data."
#    -> insert comma after "JSON string"
$d.Content.Find.Execute(
    "The above object has been converted to a JSON string which can be used to store/transmit data.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The above object has been converted to a JSON string, which can be used to store/transmit data.",
    2) | Out-Null

# 2. " Serialization formats are JSON and XML. However, XML is most used in legacy code nowadays."
#    -> "most used" becomes "mostly used"
$d.Content.Find.Execute(
    "Serialization formats are JSON and XML. However, XML is most used in legacy code nowadays.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Serialization formats are JSON and XML. However, XML is mostly used in legacy code nowadays.",
    2) | Out-Null

# 3. "The reverse process of Deserialization whereby a string is "
#    -> insert comma after "Deserialization"
$d.Content.Find.Execute(
    "The reverse process of Deserialization whereby a string is ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The reverse process of Deserialization, whereby a string is ",
    2) | Out-Null

# 4. "2) Local Storage: In JavaScript, data is inly stored as strings. Hence, the data must be serialized."
#    -> fix typo "inly" to "only"
$d.Content.Find.Execute(
    "2) Local Storage: In JavaScript, data is inly stored as strings. Hence, the data must be serialized.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2) Local Storage: In JavaScript, data is only stored as strings. Hence, the data must be serialized.",
    2) | Out-Null

# 5. "Dates becomes strings which cannot be deserialized back into Dates directly."
#    -> "becomes" replaced by "are converted into", add comma after "strings"
#       (also removes the gramStart/gramEnd proofErr markers that wrapped "becomes")
$d.Content.Find.Execute(
    "Dates becomes strings which cannot be deserialized back into Dates directly.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dates are converted into strings, which cannot be deserialized back into Dates directly.",
    2) | Out-Null

# 6. "Suppose a human wants to be teleported from point A to point B. however, the person cannot be P hysically moved..."
#    -> "however" capitalized to "However", "P" in "Physically" lowercased to "p"
$d.Content.Find.Execute(
    "Suppose a human wants to be teleported from point A to point B. however, the person cannot be Physically moved and instead ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Suppose a human wants to be teleported from point A to point B. However, the person cannot be physically moved and instead ",
    2) | Out-Null

# 7. "First, the person is scanned. Their attributes such as height, weight, memories, etc. are converted to a digital blueprint. "
#    -> add commas around "such as height, weight, memories, etc."
$d.Content.Find.Execute(
    "First, the person is scanned. Their attributes such as height, weight, memories, etc. are converted to a digital blueprint. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "First, the person is scanned. Their attributes, such as height, weight, memories, etc., are converted to a digital blueprint. ",
    2) | Out-Null

# 8. " how human transportation may not be able scan the soul of the human..."
#    -> insert "to " before "scan"
$d.Content.Find.Execute(
    "how human transportation may not be able scan the soul of the human or reconstruct emotions properly, JSON strings do not support functions and symbols, and dates cannot be converted back into Date objects from strings.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "how human transportation may not be able to scan the soul of the human or reconstruct emotions properly, JSON strings do not support functions and symbols, and dates cannot be converted back into Date objects from strings.",
    2) | Out-Null
